{"js": "// Remove the entire paragraph that documents the \"idContrato\" field under\n// the PREMIUM table description:\n//   \"idContrato: Campo que armazena o id do tipo de contrato do cliente\"\n// The paragraph (tab + text run) is deleted outright, so the surrounding\n// paragraphs (\"idUser: ...\" and the blank bold separator paragraph) become\n// adjacent.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targetText = \"idContrato: Campo que armazena o id do tipo de contrato do cliente\";\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  // The paragraph text begins with a tab character before the label, so\n  // compare against the trimmed text to be resilient to that leading tab.\n  if (para.text.trim() === targetText) {\n    para.delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# Remove the entire paragraph that documents the \"idContrato\" field under\n# the PREMIUM table description:\n#   \"idContrato: Campo que armazena o id do tipo de contrato do cliente\"\n# The paragraph (tab + text run, plus its paragraph mark) is deleted\n# outright, so the surrounding paragraphs (\"idUser: ...\" and the blank\n# bold separator paragraph) become adjacent.\n\n$d = $word.ActiveDocument\n$targetText = \"idContrato: Campo que armazena o id do tipo de contrato do cliente\"\n\nwhile ($true) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.Text = $targetText\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $found = $find.Execute()\n    if (-not $found) { break }\n\n    # Grow the found range to cover the full paragraph (the leading tab\n    # character and the trailing paragraph mark), then delete it so the\n    # paragraph disappears entirely rather than leaving an empty one behind.\n    $range.Expand(4)  # wdParagraph\n    $range.Delete()\n}\n"}
